# Update Excel SCD0011 until SCD0016
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0172 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID (column B) for both test-step rows from DGS-187 to SCD0011-003
$ws.Range("B2").Value = "SCD0011-003"
$ws.Range("B3").Value = "SCD0011-003"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6

# Scroll back to show column A and move the active selection to B4
$ws.Range("A1").Select() | Out-Null
$ws.Range("B4").Select() | Out-Null
